# Update the "Pais" (countries) COVID dashboard sheet with the latest
# refreshed numbers and re-sort consequences (Portugal overtakes Etiopia,
# Dinamarca overtakes Camerun), plus the refreshed "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: refreshed timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 15:06"

# --- Estados Unidos (row 4) --------------------------------------------------------
$ws.Range("B4").Value = 6711345
$ws.Range("C4").Value = 2887
$ws.Range("D4").Value = 3975154
$ws.Range("E4").Value = 2537629
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 198562

# --- Kuwait (row 38) ----------------------------------------------------------------
$ws.Range("B38").Value = 95472
$ws.Range("C38").Value = 708
$ws.Range("D38").Value = 85501
$ws.Range("E38").Value = 9408
$ws.Range("G38").Value = 3
$ws.Range("H38").Value = 563

# --- Paises Bajos (row 44) -----------------------------------------------------------
$ws.Range("B44").Value = 83399
$ws.Range("C44").Value = 1300
$ws.Range("G44").Value = 2
$ws.Range("H44").Value = 6256

# --- Portugal overtakes Etiopia: Portugal gets fresh numbers and moves to row 51,
#     Etiopia (unchanged numbers) drops to row 52 -------------------------------------
$ws.Range("A51").Value = "Portugal"
$ws.Range("B51").Value = 64596
$ws.Range("C51").Value = 613
$ws.Range("D51").Value = 44185
$ws.Range("E51").Value = 18540
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 1871

$ws.Range("A52").Value = "Etiopia"
$ws.Range("B52").Value = 64301
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 24983
$ws.Range("E52").Value = 38305
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 1013

# --- Dinamarca overtakes Camerun: Dinamarca gets fresh numbers and moves to row 81,
#     Camerun (unchanged numbers) drops to row 82 -------------------------------------
$ws.Range("A81").Value = "Dinamarca"
$ws.Range("B81").Value = 20237
$ws.Range("C81").Value = 347
$ws.Range("D81").Value = 16437
$ws.Range("E81").Value = 3167
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 633

$ws.Range("A82").Value = "Camerun"
$ws.Range("B82").Value = 20167
$ws.Range("C82").Value = 0
$ws.Range("D82").Value = 18837
$ws.Range("E82").Value = 915
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 415

# --- Tayikistan (row 101) ------------------------------------------------------------
$ws.Range("B101").Value = 9088
$ws.Range("C101").Value = 39
$ws.Range("D101").Value = 7853
$ws.Range("E101").Value = 1163

# --- Sri Lanka (row 135) --------------------------------------------------------------
$ws.Range("B135").Value = 3235
$ws.Range("C135").Value = 1
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = 13

# --- Timor Oriental (row 203) ---------------------------------------------------------
$ws.Range("D203").Value = 26
$ws.Range("E203").Value = 1
